$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.219.05"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "1.854.73"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'0.7031"
$ws.Range("E5").Value = "  +1.67%  "
$ws.Range("D6").Value = "'237.59"
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "'0.07995"
$ws.Range("E8").Value = "  +3.72%  "
$ws.Range("D9").Value = "'0.3021"
$ws.Range("E9").Value = "  -0.94%  "
$ws.Range("D10").Value = "'23.57"
$ws.Range("E10").Value = "  +1.34%  "
$ws.Range("D11").Value = "'0.08193"
$ws.Range("E11").Value = "  +1.63%  "
$ws.Range("D12").Value = "1.847.79"
$ws.Range("E12").Value = "  -1.82%  "
$ws.Range("D13").Value = "'5.187"
$ws.Range("E13").Value = "  -0.46%  "
$ws.Range("D14").Value = "'0.7051"
$ws.Range("E14").Value = "  -2.49%  "
$ws.Range("D15").Value = "'89.54"
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("D16").Value = "29.207.61"
$ws.Range("E16").Value = "  +0.15%  "
$ws.Range("D17").Value = "'5.819"
$ws.Range("E17").Value = "  +1.22%  "
$ws.Range("D18").Value = "'0.000007837"
$ws.Range("E18").Value = "  +0.36%  "
$ws.Range("D19").Value = "'13.19"
$ws.Range("E19").Value = "  -0.45%  "
$ws.Range("D20").Value = "'236.43"
$ws.Range("E20").Value = "  +0.79%  "
$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "  +0.14%  "
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "2.088.62"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("B23").Value = "BinanceUSD"
$ws.Range("C23").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "  +0.20%  "
$ws.Range("D24").Value = "'7.497"
$ws.Range("E24").Value = "  +0.93%  "
$ws.Range("D25").Value = "'162.53"
$ws.Range("E25").Value = "  +0.54%  "
$ws.Range("D26").Value = "'8.864"
$ws.Range("E26").Value = "  -1.18%  "
$ws.Range("D27").Value = "'0.1415"
$ws.Range("E27").Value = "  -1.38%  "
$ws.Range("D28").Value = "'18.08"
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("D29").Value = "'1.913"
$ws.Range("E29").Value = "  -2.49%  "
$ws.Range("D30").Value = "'1.401"
$ws.Range("E30").Value = "  -0.24%  "
$ws.Range("D31").Value = "'1.470"
$ws.Range("E31").Value = "  -1.28%  "
$ws.Range("D32").Value = "'4.334"
$ws.Range("E32").Value = "  -3.97%  "
$ws.Range("D33").Value = "'4.013"
$ws.Range("E33").Value = "  -0.28%  "
$ws.Range("D34").Value = "'0.05165"
$ws.Range("E34").Value = "  -0.31%  "
$ws.Range("D35").Value = "'1.166"
$ws.Range("E35").Value = "  -1.54%  "
$ws.Range("D36").Value = "'0.7108"
$ws.Range("E36").Value = "  +0.81%  "
$ws.Range("D37").Value = "'0.9949"
$ws.Range("E37").Value = "  -2.44%  "
$ws.Range("E38").Value = "  +0.30%  "
$ws.Range("E39").Value = "  -0.24%  "
$ws.Range("D40").Value = "'2.713"
$ws.Range("E40").Value = "  +1.16%  "
$ws.Range("D41").Value = "1.156.31"
$ws.Range("E41").Value = "  +4.71%  "
$ws.Range("D42").Value = "'0.9307"
$ws.Range("E42").Value = "  +0.27%  "
$ws.Range("D43").Value = "'5.982"
$ws.Range("E43").Value = "  +0.52%  "
$ws.Range("D44").Value = "'0.4253"
$ws.Range("E44").Value = "  -0.69%  "
$ws.Range("E45").Value = "  -0.10%  "
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("D47").Value = "'102.61"
$ws.Range("E47").Value = "  +0.82%  "
$ws.Range("D48").Value = "'0.5291"
$ws.Range("E48").Value = "  -2.83%  "
$ws.Range("D49").Value = "'1.740"
$ws.Range("E49").Value = "  -2.41%  "
$ws.Range("D50").Value = "'9.151"
$ws.Range("E50").Value = "  -0.27%  "
$ws.Range("B51").Value = "Aptos"
$ws.Range("C51").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D51").Value = "'6.972"
$ws.Range("E51").Value = "  -0.67%  "
